# The document contains six "<id>p029v_N</id>" sequences, each split across
# three separate runs (the literal "<id>", the bare id text, and the literal
# "</id>"). The edit merges each trio of runs into a single run that carries
# the "<id>" run's formatting (Courier New, color 7f6000, sz/szCs 18) and the
# full "<id>p029v_N</id>" text.
#
# Word's Find/Replace operates on the document's visible text stream, so
# searching for the full "<id>p029v_N</id>" string (which spans the three
# runs) and replacing it in place causes Word to collapse the matched range
# into a single run using the formatting of the range's first run -- exactly
# the merge the diff describes.

$d = $word.ActiveDocument

$ids = 1..6

foreach ($n in $ids) {
    $needle = "<id>p029v_$n</id>"
    $d.Content.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, $needle, 2)
}
